$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix lang_code column (E) for Arabic rows (11-19): was "eng", should be "ara" ---
$ws.Range("E11:E19").Value = "ara"

# --- Fix lang_code column (E) for French rows (20-28): was "eng", should be "fra" ---
$ws.Range("E20:E28").Value = "fra"

# --- Row 17 (ADD / Arabic): strip stray quote marks from name/descr ---
$ws.Range("B17").Value = "جميع تفاصيل ديموغرافية هي مطابقة"
$ws.Range("C17").Value = "جميع تفاصيل ديموغرافية هي مطابقة"

# --- Row 19 (SDM / Arabic): strip stray quote marks from name/descr ---
$ws.Range("B19").Value = "بعض التفاصيل الديمغرافية هي مطابقة"
$ws.Range("C19").Value = "بعض التفاصيل الديمغرافية هي مطابقة"

# --- Column widths for B (name) and C (descr) ---
$ws.Columns("B").ColumnWidth = 40.166666666666664
$ws.Columns("C").ColumnWidth = 45

# --- Selection / scroll position as left by the author ---
$null = $ws.Range("A7").Select()
$excel.ActiveWindow.ScrollRow = 7
$null = $ws.Range("C19").Select()

# --- Page setup (print settings) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
